$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Player Info" worksheet as the first sheet in the workbook
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the header styling used on the other sheets (bold, centered/top,
# thin boxed border)
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row (force text so IDs like "5920" stay text, matching the rest of
# the workbook where numeric-looking values are stored as text)
$playerInfo.Range("A2").Value = "'5920"
$playerInfo.Range("B2").Value = "Usman Qadir"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

# ---------------------------------------------------------------------------
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE and replace the full scorecard URL
#    with just the bare match code on the existing sheets
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "'4460"

$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "'4460"
